# Updates cryptos list price/volume figures (and the Quant/EOS row swap)
# for the "Updated cryptos list" GitHub Actions commit.
#
# All Price/Volume cells in this sheet are plain text (inline strings in
# the source OOXML), e.g. "1.003" or "30.069.07" -- NOT numbers. Excel's
# normal Range.Value assignment auto-coerces a numeric-looking string like
# "1.002" into a real number, which would change the cell's stored type
# (t="n" instead of text) and break the expected output. To avoid that we
# prefix numeric-looking replacement values with a leading apostrophe
# (Excel's standard "force text" marker) and then reset the cell's style
# back to Normal so no stray NumberFormat/quotePrefix styling is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.085.05' }
    @{ Cell = 'D3'; Value = '1.876.07' }
    @{ Cell = 'E3'; Value = '  -2.36%  ' }
    @{ Cell = 'D4'; Value = '1.002' }
    @{ Cell = 'E4'; Value = '  +0.17%  ' }
    @{ Cell = 'D5'; Value = '319.97' }
    @{ Cell = 'E5'; Value = '  -3.13%  ' }
    @{ Cell = 'D6'; Value = '1.001' }
    @{ Cell = 'E6'; Value = '  +0.12%  ' }
    @{ Cell = 'D7'; Value = '0.5041' }
    @{ Cell = 'E7'; Value = '  -3.55%  ' }
    @{ Cell = 'D8'; Value = '0.3963' }
    @{ Cell = 'E8'; Value = '  -3.31%  ' }
    @{ Cell = 'D9'; Value = '0.08211' }
    @{ Cell = 'E9'; Value = '  -3.70%  ' }
    @{ Cell = 'D10'; Value = '42.10' }
    @{ Cell = 'E10'; Value = '  -1.72%  ' }
    @{ Cell = 'E11'; Value = '  -3.06%  ' }
    @{ Cell = 'D12'; Value = '23.64' }
    @{ Cell = 'E12'; Value = '  +5.60%  ' }
    @{ Cell = 'D13'; Value = '1.868.21' }
    @{ Cell = 'E13'; Value = '  -2.67%  ' }
    @{ Cell = 'D14'; Value = '6.293' }
    @{ Cell = 'E14'; Value = '  -2.27%  ' }
    @{ Cell = 'D15'; Value = '7.199' }
    @{ Cell = 'E15'; Value = '  -2.99%  ' }
    @{ Cell = 'D16'; Value = '1.002' }
    @{ Cell = 'E16'; Value = '  +0.11%  ' }
    @{ Cell = 'D17'; Value = '91.95' }
    @{ Cell = 'E17'; Value = '  -3.83%  ' }
    @{ Cell = 'D18'; Value = '0.00001086' }
    @{ Cell = 'E18'; Value = '  -2.57%  ' }
    @{ Cell = 'D19'; Value = '0.06482' }
    @{ Cell = 'E19'; Value = '  -3.00%  ' }
    @{ Cell = 'D20'; Value = '18.14' }
    @{ Cell = 'E20'; Value = '  -1.54%  ' }
    @{ Cell = 'D21'; Value = '1.001' }
    @{ Cell = 'E21'; Value = '  +0.15%  ' }
    @{ Cell = 'D22'; Value = '30.081.80' }
    @{ Cell = 'E22'; Value = '  -0.02%  ' }
    @{ Cell = 'D23'; Value = '5.834' }
    @{ Cell = 'E23'; Value = '  -3.13%  ' }
    @{ Cell = 'D24'; Value = '11.12' }
    @{ Cell = 'E24'; Value = '  -1.95%  ' }
    @{ Cell = 'D25'; Value = '2.168' }
    @{ Cell = 'E25'; Value = '  -1.69%  ' }
    @{ Cell = 'D26'; Value = '2.088.17' }
    @{ Cell = 'E26'; Value = '  -2.22%  ' }
    @{ Cell = 'D27'; Value = '21.20' }
    @{ Cell = 'E27'; Value = '  +0.37%  ' }
    @{ Cell = 'D28'; Value = '160.58' }
    @{ Cell = 'E28'; Value = '  +0.32%  ' }
    @{ Cell = 'D29'; Value = '2.239' }
    @{ Cell = 'E29'; Value = '  -8.75%  ' }
    @{ Cell = 'D30'; Value = '127.43' }
    @{ Cell = 'E30'; Value = '  -1.52%  ' }
    @{ Cell = 'D31'; Value = '1.083' }
    @{ Cell = 'E31'; Value = '  +0.01%  ' }
    @{ Cell = 'E32'; Value = '  -1.96%  ' }
    @{ Cell = 'D33'; Value = '5.951' }
    @{ Cell = 'E33'; Value = '  -1.90%  ' }
    @{ Cell = 'D34'; Value = '3.715' }
    @{ Cell = 'E34'; Value = '  +2.36%  ' }
    @{ Cell = 'D35'; Value = '0.02430' }
    @{ Cell = 'E35'; Value = '  -2.40%  ' }
    @{ Cell = 'D36'; Value = '5.265' }
    @{ Cell = 'E36'; Value = '  +1.36%  ' }
    @{ Cell = 'E37'; Value = '  -3.98%  ' }
    @{ Cell = 'D38'; Value = '0.2138' }
    @{ Cell = 'E38'; Value = '  -3.69%  ' }
    @{ Cell = 'D39'; Value = '1.172' }
    @{ Cell = 'E39'; Value = '  -5.32%  ' }
    @{ Cell = 'D40'; Value = '8.509' }
    @{ Cell = 'E40'; Value = '  -4.34%  ' }
    @{ Cell = 'D41'; Value = '0.6293' }
    @{ Cell = 'E41'; Value = '  -3.90%  ' }
    @{ Cell = 'D42'; Value = '1.212' }
    @{ Cell = 'E42'; Value = '  -2.76%  ' }
    @{ Cell = 'E43'; Value = '  -3.33%  ' }
    @{ Cell = 'D44'; Value = '13.10' }
    @{ Cell = 'E44'; Value = '  -1.15%  ' }
    @{ Cell = 'D45'; Value = '0.5911' }
    @{ Cell = 'E45'; Value = '  -4.11%  ' }
    @{ Cell = 'D46'; Value = '2.094' }
    @{ Cell = 'E46'; Value = '  +0.39%  ' }
    @{ Cell = 'D47'; Value = '3.628' }
    @{ Cell = 'E47'; Value = '  -3.79%  ' }
    @{ Cell = 'B48'; Value = 'EOS' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos' }
    @{ Cell = 'D48'; Value = '1.208' }
    @{ Cell = 'E48'; Value = '  -3.11%  ' }
    @{ Cell = 'B49'; Value = 'Quant' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Cell = 'D49'; Value = '122.16' }
    @{ Cell = 'E49'; Value = '  -2.01%  ' }
    @{ Cell = 'D50'; Value = '77.49' }
    @{ Cell = 'E50'; Value = '  -2.95%  ' }
    @{ Cell = 'D51'; Value = '1.114' }
    @{ Cell = 'E51'; Value = '  -4.77%  ' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $value = $update.Value

    # Decide whether Excel would otherwise auto-convert this text value
    # into a number (a plain signed int/decimal like "1.002" or "42.10").
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        # Leading apostrophe forces Excel to keep it as literal text.
        $range.Value = "'" + $value
        # Drop the quote-prefix styling the assignment just introduced so
        # the cell's style index is left exactly as it was before.
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
